# Scen_PWR_BECCS_START_2027.xlsx — "UPD" sheet
#
# Commit: "Increased cost of H2, reduce ILED of BECCS to 0"
#
# The sheet is a VEDA/TIMES ~TFM_UPD transform table. Row 4 already sets
# NCAP_START = 2027 for process P-RNW-ST-BIO-CCS05 (commodity PWRBIO).
# This edit fills in the previously-blank row 5 with a second update row
# for the same process that drives NCAP_ILED (investment lead time) down
# to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attribute row: NCAP_ILED = 0 for the same process as row 4.
$ws.Range("D5").Value = "NCAP_ILED"
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = "P-RNW-ST-BIO-CCS05"

# J5 now carries real data (mirrors J4's role), so it picks up J4's
# borderless look instead of the blank hairline-bordered placeholder
# style that rows 6/7 still use.
$ws.Range("J5").Borders.LineStyle = -4142

# Leave the cursor where the author last clicked while filling the row.
$ws.Range("M5").Select() | Out-Null
